# Add checks for trigger escalability: append rows 93-100 of monitoring
# data (Date / Hour / CPU Utilization / Network In / Network Out /
# Lifecycle State) below the existing data on Sheet1, extending the used
# range from A1:H92 to A1:H100.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row, Date, Hour, CPU Utilization, Network In, Network Out, Lifecycle State
$newRows = @(
    @(93,  "2021-04-06", "03:08:44", 2.1311, "3572.0", "3204.0", "InService"),
    @(94,  "2021-04-06", "03:10:38", 2.1311, "3416.0", "3028.0", "InService"),
    @(95,  "2021-04-06", "03:11:15", 2.1311, "3416.0", "3028.0", "InService"),
    @(96,  "2021-04-06", "03:12:53", 2.3729, "3638.0", "3204.0", "InService"),
    @(97,  "2021-04-06", "03:13:18", 2.1667, "7042.0", "6843.0", "InService"),
    @(98,  "2021-04-06", "03:13:39", 2.1667, "7042.0", "6843.0", "InService"),
    @(99,  "2021-04-06", "03:14:00", 2.5,    "7042.0", "6843.0", "InService"),
    @(100, "2021-04-06", "03:14:21", 2.1667, "7042.0", "6843.0", "InService")
)

foreach ($entry in $newRows) {
    $r = $entry[0]

    # Columns A, B, D, E hold text that looks numeric/date-like (e.g. "2021-04-06",
    # "3572.0") and must be stored as literal text rather than being
    # auto-converted into a date serial / plain number. Force the cells to a
    # text format before assigning, then drop the formatting again so the
    # saved cell carries no explicit style (matching the rest of the sheet).
    $textRangeAB = $ws.Range("A$r`:B$r")
    $textRangeAB.NumberFormat = "@"
    $ws.Range("A$r").Value = $entry[1]
    $ws.Range("B$r").Value = $entry[2]
    $textRangeAB.ClearFormats()

    # Column C is a genuine number.
    $ws.Range("C$r").Value = $entry[3]

    $textRangeDE = $ws.Range("D$r`:E$r")
    $textRangeDE.NumberFormat = "@"
    $ws.Range("D$r").Value = $entry[4]
    $ws.Range("E$r").Value = $entry[5]
    $textRangeDE.ClearFormats()

    # Column H is the lifecycle state label.
    $ws.Range("H$r").Value = $entry[6]
}

Write-Host "Appended rows 93-100 to Sheet1"
